# GitHub Actions "cryptos list" refresh: updates Price (col D) and
# Volume(1h) (col E) for each coin row, and for rows 13/14 also the
# Coin name + Link (their ranking swapped: WrappedliquidstakedEther2.0
# and WrappedEther traded places).
#
# Price values are forced to text with a leading apostrophe so strings
# such as "63.00", "0.5600" or "25.859.34" keep their exact original
# formatting instead of being auto-coerced into floating point numbers
# (which would drop trailing zeros / misparse the multi-dot figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.859.34"
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = "'1.637.34"
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = "'215.34"
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").Value = "'0.5054"
$ws.Range("E6").Value = '  -0.28%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").Value = "'0.06424"
$ws.Range("E9").Value = '  +0.95%  '

$ws.Range("D10").Value = "'19.91"
$ws.Range("E10").Value = '  +1.47%  '

$ws.Range("D11").Value = "'0.07795"
$ws.Range("E11").Value = '  +0.52%  '

$ws.Range("D12").Value = "'4.286"
$ws.Range("E12").Value = '  +0.54%  '

$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.637.62"
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'1.863.51"
$ws.Range("E14").Value = '  +0.13%  '

$ws.Range("D15").Value = "'0.5600"
$ws.Range("E15").Value = '  +1.30%  '

$ws.Range("D16").Value = "'0.0₅7623"
$ws.Range("E16").Value = '  -1.40%  '

$ws.Range("D17").Value = "'63.00"
$ws.Range("E17").Value = '  -1.80%  '

$ws.Range("D18").Value = "'25.875.22"
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = "'195.08"
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("D21").Value = "'4.330"
$ws.Range("E21").Value = '  -2.54%  '

$ws.Range("D22").Value = "'9.883"
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").Value = "'6.106"
$ws.Range("E23").Value = '  +0.72%  '

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").Value = "'1.778"
$ws.Range("E25").Value = '  -6.78%  '

$ws.Range("D26").Value = "'140.39"
$ws.Range("E26").Value = '  -1.71%  '

$ws.Range("D27").Value = "'0.1257"
$ws.Range("E27").Value = '  +1.27%  '

$ws.Range("D28").Value = "'6.831"
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("D30").Value = "'1.241"
$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("D31").Value = "'0.04902"
$ws.Range("E31").Value = '  +0.68%  '

$ws.Range("E32").Value = '  +1.30%  '

$ws.Range("D33").Value = "'3.227"
$ws.Range("E33").Value = '  +1.02%  '

$ws.Range("E34").Value = '  +1.31%  '

$ws.Range("D35").Value = "'2.380"
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("D36").Value = "'0.9035"
$ws.Range("E36").Value = '  -0.29%  '

$ws.Range("D37").Value = "'2.576"
$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("D38").Value = "'0.5524"
$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("D39").Value = "'1.127.59"
$ws.Range("E39").Value = '  +0.45%  '

$ws.Range("D40").Value = "'0.01561"
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").Value = "'0.9970"

$ws.Range("D42").Value = "'5.541"
$ws.Range("E42").Value = '  -0.68%  '

$ws.Range("D43").Value = "'0.7998"
$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("D44").Value = "'98.01"
$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("D45").Value = "'1.774.07"

$ws.Range("E46").Value = '  -5.12%  '

$ws.Range("D47").Value = "'55.43"
$ws.Range("E47").Value = '  +0.84%  '

$ws.Range("D48").Value = "'0.4265"
$ws.Range("E48").Value = '  -4.32%  '

$ws.Range("D49").Value = "'7.716"
$ws.Range("E49").Value = '  +2.05%  '

$ws.Range("D50").Value = "'0.05039"
$ws.Range("E50").Value = '  -2.21%  '

$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = '  +0.40%  '
